$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H (copy formatting from the adjacent header cell, then set text)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Dominante"

# New data for column H (rows 2-4)
$ws.Range("H2").Value = "Derecha "
$ws.Range("H3").Value = "Izquierda "
$ws.Range("H4").Value = "Ambidiestro"

# Adjust column widths to fit new content (mirrors Excel's "best fit" autosize)
$ws.Columns.Item(5).ColumnWidth = 17.74
$ws.Columns.Item(8).ColumnWidth = 9.92

# Update selection to reflect post-edit active cell
$ws.Range("H5").Select() | Out-Null
